# "Chuc nang chinh ta gan hoan thien" - fill in the missing spell-check /
# grading column (H) for rows 11-15 on the "Ghi cong" sheet, and move the
# on-screen selection to where the user left off (I15, scrolled so row 4
# is at the top of the window).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ghi cong")
$ws.Activate()

# Column H was blank for these rows; the column now carries an explicit
# pass/fail (0/1) value for each person. Row 14 is the only one that picks
# up a point here, which bumps its running total in column E (a shared
# formula "=100 + SUM(F:O)") from 102 to 103.
$ws.Range("H11").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 1
$ws.Range("H15").Value = 0

# Restore the window scroll position (top-left visible cell row 4) and
# move the active selection to I15.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("I15").Select()
